# Weekly price-sheet update: insert a new week's record as row 36,
# pushing all the existing historical rows (old 36..115) down by one
# (new 37..116).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:36").Insert()

$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 45070
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112040
$ws.Range("G36").Value = "Cilantro"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = 1250
$ws.Range("N36").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 625
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = "Hortaliza"
